$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'40.818.76"
$ws.Range('D3').Value = "'2.373.63"
$ws.Range('E3').Value = '  -3.99%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = "'312.68"
$ws.Range('E5').Value = '  -1.81%  '
$ws.Range('D6').Value = "'88.00"
$ws.Range('E6').Value = '  -5.69%  '
$ws.Range('D7').Value = "'0.530"
$ws.Range('E7').Value = '  -3.90%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -4.53%  '
$ws.Range('D10').Value = "'0.0837"
$ws.Range('E10').Value = '  -3.69%  '
$ws.Range('D11').Value = "'30.93"
$ws.Range('E11').Value = '  -7.51%  '
$ws.Range('E12').Value = '  -1.39%  '
$ws.Range('D13').Value = "'2.738.56"
$ws.Range('E13').Value = '  -4.05%  '
$ws.Range('D14').Value = "'6.56"
$ws.Range('E14').Value = '  -5.13%  '
$ws.Range('D15').Value = "'15.03"
$ws.Range('E15').Value = '  -4.35%  '
$ws.Range('D16').Value = "'2.356.23"
$ws.Range('E16').Value = '  -4.46%  '
$ws.Range('D17').Value = "'0.761"
$ws.Range('E17').Value = '  -3.97%  '
$ws.Range('D18').Value = "'40.698.64"
$ws.Range('E18').Value = '  -2.35%  '
$ws.Range('D19').Value = "'0.0₃0911"
$ws.Range('E19').Value = '  -4.01%  '
$ws.Range('D20').Value = "'6.16"
$ws.Range('E20').Value = '  -4.83%  '
$ws.Range('D21').Value = "'69.39"
$ws.Range('E21').Value = '  -2.38%  '
$ws.Range('D22').Value = "'10.77"
$ws.Range('E22').Value = '  -4.60%  '
$ws.Range('D23').Value = "'232.58"
$ws.Range('E23').Value = '  -2.85%  '
$ws.Range('D24').Value = "'2.65"
$ws.Range('E24').Value = '  -3.44%  '
$ws.Range('E25').Value = '  +0.15%  '
$ws.Range('D26').Value = "'1.81"
$ws.Range('E26').Value = '  -6.60%  '
$ws.Range('D27').Value = "'23.77"
$ws.Range('E27').Value = '  -3.75%  '
$ws.Range('D28').Value = "'2.17"
$ws.Range('E28').Value = '  -3.92%  '
$ws.Range('D29').Value = "'9.32"
$ws.Range('E29').Value = '  -4.69%  '
$ws.Range('D30').Value = "'33.69"
$ws.Range('E30').Value = '  -6.73%  '
$ws.Range('D31').Value = "'153.62"
$ws.Range('E31').Value = '  -4.16%  '
$ws.Range('E32').Value = '  +0.00%  '
$ws.Range('E33').Value = '  -5.39%  '
$ws.Range('D34').Value = "'0.0732"
$ws.Range('E34').Value = '  -4.30%  '
$ws.Range('D35').Value = "'2.44"
$ws.Range('E35').Value = '  -5.49%  '
$ws.Range('E36').Value = '  -2.24%  '
$ws.Range('D37').Value = "'2.78"
$ws.Range('E37').Value = '  -4.86%  '
$ws.Range('D38').Value = "'15.87"
$ws.Range('E38').Value = '  -9.77%  '
$ws.Range('D39').Value = "'0.0999"
$ws.Range('E39').Value = '  -4.04%  '
$ws.Range('D40').Value = "'1.72"
$ws.Range('E40').Value = '  -7.88%  '
$ws.Range('D41').Value = "'3.82"
$ws.Range('E41').Value = '  -4.68%  '
$ws.Range('D42').Value = "'2.36"
$ws.Range('E42').Value = '  -6.58%  '
$ws.Range('D43').Value = "'1.953.26"
$ws.Range('E43').Value = '  -2.27%  '
$ws.Range('D44').Value = "'0.0270"
$ws.Range('E44').Value = '  -5.19%  '
$ws.Range('D45').Value = "'17.47"
$ws.Range('E45').Value = '  -7.38%  '
$ws.Range('D46').Value = "'2.76"
$ws.Range('E46').Value = '  -7.98%  '
$ws.Range('E47').Value = '  -1.76%  '
$ws.Range('D48').Value = "'2.603.02"
$ws.Range('E48').Value = '  -3.97%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').Value = "'93.55"
$ws.Range('E49').Value = '  -4.26%  '
$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D50').Value = "'72.76"
$ws.Range('E50').Value = '  -1.47%  '
$ws.Range('D51').Value = "'50.63"
$ws.Range('E51').Value = '  -3.78%  '
